$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.9818413333333332
$ws.Range("H2").Value = 2.945524
$ws.Range("I2").Value = 0.0542081715590086
$ws.Range("J2").Value = 0.05420817155900861
$ws.Range("M2").Value = 0.003710666666666666
$ws.Range("N2").Value = 0.011132
$ws.Range("O2").Value = 0.001642024256586498
$ws.Range("P2").Value = 0.001642024256586498
$ws.Range("Q2").Value = 0.003643285907555555
$ws.Range("R2").Value = 0.032789573168
$ws.Range("S2").Value = 0.00008901113260509446
$ws.Range("T2").Value = 0.00008901113260509447
$ws.Range("G3").Value = 0.9818413333333332
$ws.Range("H3").Value = 2.945524
$ws.Range("I3").Value = 0.0542081715590086
$ws.Range("J3").Value = 0.05420817155900861
$ws.Range("M3").Value = 2.049608666666666
$ws.Range("N3").Value = 6.148826
$ws.Range("O3").Value = 0.9069818039462568
$ws.Range("P3").Value = 0.9069818039462569
$ws.Range("Q3").Value = 2.012390506091555
$ws.Range("R3").Value = 18.111514554824
$ws.Range("S3").Value = 0.04916582522921779
$ws.Range("T3").Value = 0.04916582522921781
$ws.Range("G4").Value = 0.9818413333333332
$ws.Range("H4").Value = 2.945524
$ws.Range("I4").Value = 0.0542081715590086
$ws.Range("J4").Value = 0.05420817155900861
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.206493
$ws.Range("N4").Value = 0.619479
$ws.Range("O4").Value = 0.09137617179715662
$ws.Range("P4").Value = 0.09137617179715661
$ws.Range("Q4").Value = 0.202743362444
$ws.Range("R4").Value = 1.824690261996
$ws.Range("S4").Value = 0.00495333519718571
$ws.Range("T4").Value = 0.004953335197185709
$ws.Range("G5").Value = 9.968049999999998
$ws.Range("H5").Value = 29.90414999999999
$ws.Range("I5").Value = 0.5503432643992467
$ws.Range("J5").Value = 0.5503432643992469
$ws.Range("M5").Value = 0.003710666666666666
$ws.Range("N5").Value = 0.011132
$ws.Range("O5").Value = 0.001642024256586498
$ws.Range("P5").Value = 0.001642024256586498
$ws.Range("Q5").Value = 0.03698811086666666
$ws.Range("R5").Value = 0.3328929977999999
$ws.Range("S5").Value = 0.0009036769895925599
$ws.Range("T5").Value = 0.00090367698959256
$ws.Range("G6").Value = 9.968049999999998
$ws.Range("H6").Value = 29.90414999999999
$ws.Range("I6").Value = 0.5503432643992467
$ws.Range("J6").Value = 0.5503432643992469
$ws.Range("M6").Value = 2.049608666666666
$ws.Range("N6").Value = 6.148826
$ws.Range("O6").Value = 0.9069818039462568
$ws.Range("P6").Value = 0.9069818039462569
$ws.Range("Q6").Value = 20.43060166976666
$ws.Range("R6").Value = 183.8754150278999
$ws.Range("S6").Value = 0.4991513267345006
$ws.Range("T6").Value = 0.4991513267345007
$ws.Range("G7").Value = 9.968049999999998
$ws.Range("H7").Value = 29.90414999999999
$ws.Range("I7").Value = 0.5503432643992467
$ws.Range("J7").Value = 0.5503432643992469
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.206493
$ws.Range("N7").Value = 0.619479
$ws.Range("O7").Value = 0.09137617179715662
$ws.Range("P7").Value = 0.09137617179715661
$ws.Range("Q7").Value = 2.05833254865
$ws.Range("R7").Value = 18.52499293785
$ws.Range("S7").Value = 0.05028826067515356
$ws.Range("T7").Value = 0.05028826067515357
$ws.Range("G8").Value = 7.162531666666666
$ws.Range("H8").Value = 21.487595
$ws.Range("I8").Value = 0.3954485640417445
$ws.Range("J8").Value = 0.3954485640417446
$ws.Range("M8").Value = 0.003710666666666666
$ws.Range("N8").Value = 0.011132
$ws.Range("O8").Value = 0.001642024256586498
$ws.Range("P8").Value = 0.001642024256586498
$ws.Range("Q8").Value = 0.02657776750444444
$ws.Range("R8").Value = 0.23919990754
$ws.Range("S8").Value = 0.0006493361343888438
$ws.Range("T8").Value = 0.0006493361343888438
$ws.Range("G9").Value = 7.162531666666666
$ws.Range("H9").Value = 21.487595
$ws.Range("I9").Value = 0.3954485640417445
$ws.Range("J9").Value = 0.3954485640417446
$ws.Range("M9").Value = 2.049608666666666
$ws.Range("N9").Value = 6.148826
$ws.Range("O9").Value = 0.9069818039462568
$ws.Range("P9").Value = 0.9069818039462569
$ws.Range("Q9").Value = 14.68038697927444
$ws.Range("R9").Value = 132.12348281347
$ws.Range("S9").Value = 0.3586646519825383
$ws.Range("T9").Value = 0.3586646519825384
$ws.Range("G10").Value = 7.162531666666666
$ws.Range("H10").Value = 21.487595
$ws.Range("I10").Value = 0.3954485640417445
$ws.Range("J10").Value = 0.3954485640417446
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.206493
$ws.Range("N10").Value = 0.619479
$ws.Range("O10").Value = 0.09137617179715662
$ws.Range("P10").Value = 0.09137617179715661
$ws.Range("Q10").Value = 1.479012651445
$ws.Range("R10").Value = 13.311113863005
$ws.Range("S10").Value = 0.03613457592481734
$ws.Range("T10").Value = 0.03613457592481734
